$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at 535, pushing existing rows 535-545 down to 542-552
$ws.Rows("535:541").Insert()

# Row 535
$ws.Cells.Item(535, 1).Value = 6
$ws.Cells.Item(535, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(535, 3).Value = "Metropolitana"
$ws.Cells.Item(535, 4).Value = 45239
$ws.Cells.Item(535, 5).Value = 13
$ws.Cells.Item(535, 6).Value = "Fruta"
$ws.Cells.Item(535, 7).Value = 100107
$ws.Cells.Item(535, 8).Value = "Otros"
$ws.Cells.Item(535, 9).Value = 100107002
$ws.Cells.Item(535, 10).Value = "Chirimoya"
$ws.Cells.Item(535, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(535, 12).Value = "Especial"
$ws.Cells.Item(535, 13).Value = 275
$ws.Cells.Item(535, 14).Value = 19000
$ws.Cells.Item(535, 15).Value = 19000
$ws.Cells.Item(535, 16).Value = 19000
$ws.Cells.Item(535, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(535, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(535, 19).Value = 1900
$ws.Cells.Item(535, 20).Value = 10

# Row 536
$ws.Cells.Item(536, 1).Value = 6
$ws.Cells.Item(536, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 45239
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = "Fruta"
$ws.Cells.Item(536, 7).Value = 100107
$ws.Cells.Item(536, 8).Value = "Otros"
$ws.Cells.Item(536, 9).Value = 100107002
$ws.Cells.Item(536, 10).Value = "Chirimoya"
$ws.Cells.Item(536, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(536, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(536, 13).Value = 180
$ws.Cells.Item(536, 14).Value = 18000
$ws.Cells.Item(536, 15).Value = 18000
$ws.Cells.Item(536, 16).Value = 18000
$ws.Cells.Item(536, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(536, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(536, 19).Value = 1800
$ws.Cells.Item(536, 20).Value = 10

# Row 537
$ws.Cells.Item(537, 1).Value = 6
$ws.Cells.Item(537, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(537, 3).Value = "Metropolitana"
$ws.Cells.Item(537, 4).Value = 45239
$ws.Cells.Item(537, 5).Value = 13
$ws.Cells.Item(537, 6).Value = "Fruta"
$ws.Cells.Item(537, 7).Value = 100107
$ws.Cells.Item(537, 8).Value = "Otros"
$ws.Cells.Item(537, 9).Value = 100107002
$ws.Cells.Item(537, 10).Value = "Chirimoya"
$ws.Cells.Item(537, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(537, 12).Value = "Primera"
$ws.Cells.Item(537, 13).Value = 400
$ws.Cells.Item(537, 14).Value = 15000
$ws.Cells.Item(537, 15).Value = 15000
$ws.Cells.Item(537, 16).Value = 15000
$ws.Cells.Item(537, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(537, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(537, 19).Value = 1500
$ws.Cells.Item(537, 20).Value = 10

# Row 538
$ws.Cells.Item(538, 1).Value = 6
$ws.Cells.Item(538, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(538, 3).Value = "Metropolitana"
$ws.Cells.Item(538, 4).Value = 45239
$ws.Cells.Item(538, 5).Value = 13
$ws.Cells.Item(538, 6).Value = "Fruta"
$ws.Cells.Item(538, 7).Value = 100107
$ws.Cells.Item(538, 8).Value = "Otros"
$ws.Cells.Item(538, 9).Value = 100107002
$ws.Cells.Item(538, 10).Value = "Chirimoya"
$ws.Cells.Item(538, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(538, 12).Value = "Primera"
$ws.Cells.Item(538, 13).Value = 150
$ws.Cells.Item(538, 14).Value = 15000
$ws.Cells.Item(538, 15).Value = 15000
$ws.Cells.Item(538, 16).Value = 15000
$ws.Cells.Item(538, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(538, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(538, 19).Value = 1500
$ws.Cells.Item(538, 20).Value = 10

# Row 539
$ws.Cells.Item(539, 1).Value = 6
$ws.Cells.Item(539, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(539, 3).Value = "Metropolitana"
$ws.Cells.Item(539, 4).Value = 45239
$ws.Cells.Item(539, 5).Value = 13
$ws.Cells.Item(539, 6).Value = "Fruta"
$ws.Cells.Item(539, 7).Value = 100107
$ws.Cells.Item(539, 8).Value = "Otros"
$ws.Cells.Item(539, 9).Value = 100107002
$ws.Cells.Item(539, 10).Value = "Chirimoya"
$ws.Cells.Item(539, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(539, 12).Value = "Segunda"
$ws.Cells.Item(539, 13).Value = 350
$ws.Cells.Item(539, 14).Value = 13000
$ws.Cells.Item(539, 15).Value = 13000
$ws.Cells.Item(539, 16).Value = 13000
$ws.Cells.Item(539, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(539, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(539, 19).Value = 1300
$ws.Cells.Item(539, 20).Value = 10

# Row 540
$ws.Cells.Item(540, 1).Value = 6
$ws.Cells.Item(540, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(540, 3).Value = "Metropolitana"
$ws.Cells.Item(540, 4).Value = 45239
$ws.Cells.Item(540, 5).Value = 13
$ws.Cells.Item(540, 6).Value = "Fruta"
$ws.Cells.Item(540, 7).Value = 100107
$ws.Cells.Item(540, 8).Value = "Otros"
$ws.Cells.Item(540, 9).Value = 100107002
$ws.Cells.Item(540, 10).Value = "Chirimoya"
$ws.Cells.Item(540, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(540, 12).Value = "Segunda"
$ws.Cells.Item(540, 13).Value = 140
$ws.Cells.Item(540, 14).Value = 12000
$ws.Cells.Item(540, 15).Value = 12000
$ws.Cells.Item(540, 16).Value = 12000
$ws.Cells.Item(540, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(540, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(540, 19).Value = 1200
$ws.Cells.Item(540, 20).Value = 10

# Row 541
$ws.Cells.Item(541, 1).Value = 6
$ws.Cells.Item(541, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(541, 3).Value = "Metropolitana"
$ws.Cells.Item(541, 4).Value = 45239
$ws.Cells.Item(541, 5).Value = 13
$ws.Cells.Item(541, 6).Value = "Fruta"
$ws.Cells.Item(541, 7).Value = 100107
$ws.Cells.Item(541, 8).Value = "Otros"
$ws.Cells.Item(541, 9).Value = 100107002
$ws.Cells.Item(541, 10).Value = "Chirimoya"
$ws.Cells.Item(541, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(541, 12).Value = "Tercera"
$ws.Cells.Item(541, 13).Value = 120
$ws.Cells.Item(541, 14).Value = 10000
$ws.Cells.Item(541, 15).Value = 10000
$ws.Cells.Item(541, 16).Value = 10000
$ws.Cells.Item(541, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(541, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(541, 19).Value = 1000
$ws.Cells.Item(541, 20).Value = 10
